# The workbook gained one new data row (weekly price record) that was
# inserted right before the existing row 376, pushing all subsequent
# rows (376..472) down by one (to 377..473). The sheet's used range
# grows from A1:R472 to A1:R473 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 376, shifting rows 376-472 down to 377-473
# and carrying their formatting (date style on column D, etc.) with them.
$ws.Rows(376).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A376").Value = 6
$ws.Range("B376").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C376").Value = "Metropolitana"
$ws.Range("D376").Value = 44782
$ws.Range("E376").Value = 13
$ws.Range("F376").Value = 100112043
$ws.Range("G376").Value = "Pepino ensalada"
$ws.Range("H376").Value = "Sin especificar"
$ws.Range("I376").Value = "Primera"
$ws.Range("J376").Value = 250
$ws.Range("K376").Value = 19000
$ws.Range("L376").Value = 20000
$ws.Range("M376").Value = 19600
$ws.Range("N376").Value = '$/caja 60 unidades'
$ws.Range("O376").Value = "Región de Arica y Parinacota"
$ws.Range("P376").Value = 327
$ws.Range("Q376").Value = 60
$ws.Range("R376").Value = "Hortaliza"
